$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("codeforiati:group-code", "codeforiati:category-name"),
    @("110", "Education, Level Unspecified"),
    @("110", "Education, Level Unspecified"),
    @("110", "Education, Level Unspecified"),
    @("110", "Education, Level Unspecified"),
    @("110", "Basic Education"),
    @("110", "Basic Education"),
    @("110", "Basic Education"),
    @("110", "Basic Education"),
    @("110", "Basic Education"),
    @("110", "Basic Education"),
    @("110", "Basic Education"),
    @("110", "Secondary Education"),
    @("110", "Secondary Education"),
    @("110", "Post-Secondary Education"),
    @("110", "Post-Secondary Education"),
    @("120", "Health, General"),
    @("120", "Health, General"),
    @("120", "Health, General"),
    @("120", "Health, General"),
    @("120", "Basic Health"),
    @("120", "Basic Health"),
    @("120", "Basic Health"),
    @("120", "Basic Health"),
    @("120", "Basic Health"),
    @("120", "Basic Health"),
    @("120", "Basic Health"),
    @("120", "Basic Health"),
    @("120", "Basic Health"),
    @("120", "Non-communicable diseases (NCDs)"),
    @("120", "Non-communicable diseases (NCDs)"),
    @("120", "Non-communicable diseases (NCDs)"),
    @("120", "Non-communicable diseases (NCDs)"),
    @("120", "Non-communicable diseases (NCDs)"),
    @("120", "Non-communicable diseases (NCDs)"),
    @("130", "Population Policies/Programmes & Reproductive Health"),
    @("130", "Population Policies/Programmes & Reproductive Health"),
    @("130", "Population Policies/Programmes & Reproductive Health"),
    @("130", "Population Policies/Programmes & Reproductive Health"),
    @("130", "Population Policies/Programmes & Reproductive Health"),
    @("140", "Water Supply & Sanitation"),
    @("140", "Water Supply & Sanitation"),
    @("140", "Water Supply & Sanitation"),
    @("140", "Water Supply & Sanitation"),
    @("140", "Water Supply & Sanitation"),
    @("140", "Water Supply & Sanitation"),
    @("140", "Water Supply & Sanitation"),
    @("140", "Water Supply & Sanitation"),
    @("140", "Water Supply & Sanitation"),
    @("140", "Water Supply & Sanitation"),
    @("140", "Water Supply & Sanitation"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Government & Civil Society-general"),
    @("150", "Conflict, Peace & Security"),
    @("150", "Conflict, Peace & Security"),
    @("150", "Conflict, Peace & Security"),
    @("150", "Conflict, Peace & Security"),
    @("150", "Conflict, Peace & Security"),
    @("150", "Conflict, Peace & Security"),
    @("160", "Other Social Infrastructure & Services"),
    @("160", "Other Social Infrastructure & Services"),
    @("160", "Other Social Infrastructure & Services"),
    @("160", "Other Social Infrastructure & Services"),
    @("160", "Other Social Infrastructure & Services"),
    @("160", "Other Social Infrastructure & Services"),
    @("160", "Other Social Infrastructure & Services"),
    @("160", "Other Social Infrastructure & Services"),
    @("160", "Other Social Infrastructure & Services"),
    @("160", "Other Social Infrastructure & Services"),
    @("160", "Other Social Infrastructure & Services"),
    @("210", "Transport & Storage"),
    @("210", "Transport & Storage"),
    @("210", "Transport & Storage"),
    @("210", "Transport & Storage"),
    @("210", "Transport & Storage"),
    @("210", "Transport & Storage"),
    @("210", "Transport & Storage"),
    @("220", "Communications"),
    @("220", "Communications"),
    @("220", "Communications"),
    @("220", "Communications"),
    @("230", "Energy Policy"),
    @("230", "Energy Policy"),
    @("230", "Energy Policy"),
    @("230", "Energy Policy"),
    @("230", "Energy generation, renewable sources"),
    @("230", "Energy generation, renewable sources"),
    @("230", "Energy generation, renewable sources"),
    @("230", "Energy generation, renewable sources"),
    @("230", "Energy generation, renewable sources"),
    @("230", "Energy generation, renewable sources"),
    @("230", "Energy generation, renewable sources"),
    @("230", "Energy generation, renewable sources"),
    @("230", "Energy generation, renewable sources"),
    @("230", "Energy generation, non-renewable sources"),
    @("230", "Energy generation, non-renewable sources"),
    @("230", "Energy generation, non-renewable sources"),
    @("230", "Energy generation, non-renewable sources"),
    @("230", "Energy generation, non-renewable sources"),
    @("230", "Energy generation, non-renewable sources"),
    @("230", "Hybrid energy plants"),
    @("230", "Nuclear energy plants"),
    @("230", "Energy distribution"),
    @("230", "Energy distribution"),
    @("230", "Energy distribution"),
    @("230", "Energy distribution"),
    @("230", "Energy distribution"),
    @("230", "Energy distribution"),
    @("230", "Energy distribution"),
    @("240", "Banking & Financial Services"),
    @("240", "Banking & Financial Services"),
    @("240", "Banking & Financial Services"),
    @("240", "Banking & Financial Services"),
    @("240", "Banking & Financial Services"),
    @("240", "Banking & Financial Services"),
    @("250", "Business & Other Services"),
    @("250", "Business & Other Services"),
    @("250", "Business & Other Services"),
    @("250", "Business & Other Services"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Agriculture"),
    @("310", "Forestry"),
    @("310", "Forestry"),
    @("310", "Forestry"),
    @("310", "Forestry"),
    @("310", "Forestry"),
    @("310", "Forestry"),
    @("310", "Fishing"),
    @("310", "Fishing"),
    @("310", "Fishing"),
    @("310", "Fishing"),
    @("310", "Fishing"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Industry"),
    @("320", "Mineral Resources & Mining"),
    @("320", "Mineral Resources & Mining"),
    @("320", "Mineral Resources & Mining"),
    @("320", "Mineral Resources & Mining"),
    @("320", "Mineral Resources & Mining"),
    @("320", "Mineral Resources & Mining"),
    @("320", "Mineral Resources & Mining"),
    @("320", "Mineral Resources & Mining"),
    @("320", "Mineral Resources & Mining"),
    @("320", "Mineral Resources & Mining"),
    @("320", "Construction"),
    @("331", "Trade Policies & Regulations"),
    @("331", "Trade Policies & Regulations"),
    @("331", "Trade Policies & Regulations"),
    @("331", "Trade Policies & Regulations"),
    @("331", "Trade Policies & Regulations"),
    @("331", "Trade Policies & Regulations"),
    @("332", "Tourism"),
    @("410", "General Environment Protection"),
    @("410", "General Environment Protection"),
    @("410", "General Environment Protection"),
    @("410", "General Environment Protection"),
    @("410", "General Environment Protection"),
    @("410", "General Environment Protection"),
    @("430", "Other Multisector"),
    @("430", "Other Multisector"),
    @("430", "Other Multisector"),
    @("430", "Other Multisector"),
    @("430", "Other Multisector"),
    @("430", "Other Multisector"),
    @("430", "Other Multisector"),
    @("430", "Other Multisector"),
    @("430", "Other Multisector"),
    @("430", "Other Multisector"),
    @("510", "General Budget Support"),
    @("520", "Development Food Assistance"),
    @("530", "Other Commodity Assistance"),
    @("530", "Other Commodity Assistance"),
    @("600", "Action Relating to Debt"),
    @("600", "Action Relating to Debt"),
    @("600", "Action Relating to Debt"),
    @("600", "Action Relating to Debt"),
    @("600", "Action Relating to Debt"),
    @("600", "Action Relating to Debt"),
    @("600", "Action Relating to Debt"),
    @("720", "Emergency Response"),
    @("720", "Emergency Response"),
    @("720", "Emergency Response"),
    @("730", "Reconstruction Relief & Rehabilitation"),
    @("740", "Disaster Prevention & Preparedness"),
    @("910", "Administrative Costs of Donors"),
    @("930", "Refugees in Donor Countries"),
    @("998", "Unallocated / Unspecified"),
    @("998", "Unallocated / Unspecified")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 1
    $pair = $data[$i]
    $ws.Cells.Item($rowNum, 5).Value = $pair[0]
    $ws.Cells.Item($rowNum, 6).Value = $pair[1]
}

Write-Host "Done applying SectorGroup E/F column corrections"